$wb = $excel.ActiveWorkbook
$wsJ = $wb.Worksheets.Item("Jugadores")

# Player stats rows 2-16 (row 11 "Nico" is new; rows 12-16 shift down from old 11-15;
# Subra row grows from 15 to 16 with updated values; dimension becomes A1:W16)
$row = @("Alberto R", 0, 6, 3, 3, 0, 34.05, 12, 3, 7, 2, 8, 5, 15, 0, 2, 1.5, 33.05, 45.85, 25, 0, 0, 2)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(2, $c).Value = $row[$c - 1] }

$row = @("Feylin", 7, 8, 3, 1, 0, 55.5, 15, 6, 11, 0, 0, 6, 11, 3, 7, 1.25, 54.18, 54.18, 0, 37.5, 85, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(3, $c).Value = $row[$c - 1] }

$row = @("Gustavo", 3, 13, 15, 2, 0, 71.1, 12, 3, 10, 1, 6, 4, 16, 3, 7, 2.75, 21.25, 22.5, 12.5, 37.5, 11, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(4, $c).Value = $row[$c - 1] }

$row = @("Jorge", 0, 1, 2, 0, 0, 8.9, 2, 1, 3, 0, 1, 1, 4, 0, 0, 1, 25, 33.3, 0, 0, 10, 1)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(5, $c).Value = $row[$c - 1] }

$row = @("Jurque", 9, 14, 1, 4, 0, 68.37, 20, 9, 18, 0, 0, 9, 18, 2, 7, 2.5, 54.15, 54.15, 0, 25, 2, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(6, $c).Value = $row[$c - 1] }

$row = @("Leonard", 4, 11, 8, 10, 0, 76.2, 41, 10, 18, 5, 20, 15, 38, 6, 15, 3.5, 40.5, 65.9, 24.3, 48.32, 90, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(7, $c).Value = $row[$c - 1] }

$row = @("Lianes", 10, 11, 8, 11, 1, 90.52, 48, 4, 12, 12, 34, 16, 46, 4, 7, 1.75, 37.08, 28.75, 37.08, 50, 14, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(8, $c).Value = $row[$c - 1] }

$row = @("Marín", 8, 4, 5, 3, 0, 47.82, 15, 1, 5, 4, 11, 5, 16, 1, 2, 1.25, 32.5, 8.32, 35, 12.5, 87, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(9, $c).Value = $row[$c - 1] }

$row = @("Max", 13, 8, 8, 3, 6, 71.72, 30, 12, 18, 0, 3, 12, 21, 6, 9, 1.5, 61.18, 67.85, 0, 68.75, 16, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(10, $c).Value = $row[$c - 1] }

$row = @("Nico", 0, 2, 0, 0, 0, 2.8, 0, 0, 1, 0, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 42, 1)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(11, $c).Value = $row[$c - 1] }

$row = @("Oscar", 3, 0, 0, 1, 0, 12.17, 2, 1, 3, 0, 0, 1, 3, 0, 0, 0.5, 50, 50, 0, 0, 31, 2)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(12, $c).Value = $row[$c - 1] }

$row = @("Pete", 5, 7, 3, 0, 1, 32.3, 8, 4, 7, 0, 1, 4, 8, 0, 1, 3, 26.67, 26.67, 0, 0, 4, 3)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(13, $c).Value = $row[$c - 1] }

$row = @("Raúl", 0, 4, 3, 2, 0, 42.87, 43, 3, 5, 10, 24, 13, 29, 7, 10, 2.33, 44.77, 77.77, 42.5, 47.23, 7, 3)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(14, $c).Value = $row[$c - 1] }

$row = @("Rubén", 2, 6, 2, 2, 0, 41.75, 6, 2, 5, 0, 1, 2, 6, 2, 6, 2, 25, 33.35, 0, 37.5, 8, 2)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(15, $c).Value = $row[$c - 1] }

$row = @("Subra", 20, 24, 8, 3, 9, 93.72, 44, 16, 30, 0, 1, 16, 31, 12, 22, 1.25, 53.75, 55.42, 0, 63.25, 93, 4)
for ($c = 1; $c -le $row.Length; $c++) { $wsJ.Cells.Item(16, $c).Value = $row[$c - 1] }

$wsR = $wb.Worksheets.Item("Stats_Rival")
$row2 = @(74.5, 55.75, 21, 30, 51, 7.5, 24, 31.5)
for ($c = 1; $c -le $row2.Length; $c++) { $wsR.Cells.Item(2, $c).Value = $row2[$c - 1] }